# Add "Start time" and "End time" columns to the Groups sheet, between
# the existing "When" and "Venue" columns, and leave the Groups sheet as
# the active sheet/selection (matching the author's commit:
# "start and end times added - Change to allow Group Start and End times").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groups")

# Insert two new blank columns where the new headers should go (columns
# H and I), pushing "Venue" and everything after it two columns to the
# right (H:O -> J:Q).
$ws.Columns("H:I").Insert()

# Populate the new header cells.
$ws.Cells.Item(1, 8).Value = "Start time"
$ws.Cells.Item(1, 9).Value = "End time"

# Make Groups the active sheet and select the new "End time" header cell,
# mirroring the saved selection/active-tab state in the target workbook.
$ws.Activate() | Out-Null
$ws.Range("I1").Select() | Out-Null
